$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook reports localization handoff status for two files:
#   14a3b362-714d-4008-932b-2d15393f2e8c.md   (still in flight)
#   19475d89-b2f9-4755-80cd-2a4b0b5bc60b.md   (dropped from this report)
# plus the always-present ".localization-config" row.
#
# The new report run:
#   - flips the 14a3b362 entry's status from "Handed back" to
#     "Not yet handed off" (it has not come back from localization yet)
#   - refreshes the latest-handoff timestamps for that file
#   - drops the 19475d89 entry entirely (its row is removed, and every
#     sheet loses one row)
# ---------------------------------------------------------------------------

$urlOltest        = "https://github.com/OpenLocalizationTest/oltest/blob/5f9f0af19f18313ec636f9f48b72c8748c7a8544"
$urlOltestZhCn     = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ac97d6806093d5ec0a4730e727bd8d40f158177d"
$urlOltestDeDe     = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/495bfa5db353ff46001df1abedefa2f9a316ca6d"
$urlHandoffZhCn    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ffeafd63e0ededa62305c16c03a94f81ca7b04f6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang"
$urlHandbackZhCn   = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/880eb6945afc60531bd8c1b7005d6ec89e435a23/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang"
$urlHandoffDeDe    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4918f75884f48508f7094f91dbaceb7fba13ca6f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang"
$urlHandbackDeDe   = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6770de6f3959863144f0f5da977af477a095e714/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang"

$mdName    = "14a3b362-714d-4008-932b-2d15393f2e8c.md"
$zhCnXlf   = "14a3b362-714d-4008-932b-2d15393f2e8c.309b9c75701238e0f13f0b6d6ff1d0e54ec7fc20.zh-cn.xlf"
$deDeXlf   = "14a3b362-714d-4008-932b-2d15393f2e8c.309b9c75701238e0f13f0b6d6ff1d0e54ec7fc20.de-de.xlf"
$cfgName   = ".localization-config"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ws = $wb.Worksheets.Item("Overview")

# Update status for the 14a3b362 row before the other row is removed.
$ws.Range("B2").Value = "Not yet handed off"
$ws.Range("C2").Value = "Not yet handed off"

# Drop the 19475d89 row (row 3) entirely; .localization-config shifts up to row 3.
$ws.Hyperlinks.Delete()
$ws.Rows(3).Delete()

# Re-create the hyperlinks that survive, in left-to-right / top-to-bottom order.
$ws.Hyperlinks.Add($ws.Range("A2"), "$urlOltest/e2e/$mdName", "", "", $mdName)
$ws.Hyperlinks.Add($ws.Range("A3"), "$urlOltest/$cfgName", "", "", $cfgName)

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B2").Value = "Not yet handed off"
$ws.Range("D2").Value = "2016-01-11 03:02:54"

$ws.Hyperlinks.Delete()
$ws.Rows(3).Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "$urlOltest/e2e/$mdName", "", "", $mdName)
$ws.Hyperlinks.Add($ws.Range("C2"), "$urlHandoffZhCn/$zhCnXlf", "", "", $zhCnXlf)
$ws.Hyperlinks.Add($ws.Range("E2"), "$urlOltestZhCn/e2e/$mdName", "", "", $mdName)
$ws.Hyperlinks.Add($ws.Range("F2"), "$urlHandbackZhCn/$zhCnXlf", "", "", $zhCnXlf)
$ws.Hyperlinks.Add($ws.Range("A3"), "$urlOltest/$cfgName", "", "", $cfgName)

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B2").Value = "Not yet handed off"
$ws.Range("D2").Value = "2016-01-11 03:03:10"

$ws.Hyperlinks.Delete()
$ws.Rows(3).Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "$urlOltest/e2e/$mdName", "", "", $mdName)
$ws.Hyperlinks.Add($ws.Range("C2"), "$urlHandoffDeDe/$deDeXlf", "", "", $deDeXlf)
$ws.Hyperlinks.Add($ws.Range("E2"), "$urlOltestDeDe/e2e/$mdName", "", "", $mdName)
$ws.Hyperlinks.Add($ws.Range("F2"), "$urlHandbackDeDe/$deDeXlf", "", "", $deDeXlf)
$ws.Hyperlinks.Add($ws.Range("A3"), "$urlOltest/$cfgName", "", "", $cfgName)
